$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the BRACKET support prices in column D (rows 33-37)
$ws.Range("D33").Value = 767.647
$ws.Range("D34").Value = 1139.001
$ws.Range("D35").Value = 1427.198
$ws.Range("D36").Value = 1718.204
$ws.Range("D37").Value = 1878.96
